$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Story Points for rows 8 and 9 (C8, C9): 3 -> 1.5
$ws.Range("C8").Value = 1.5
$ws.Range("C9").Value = 1.5

# Add Story Points value for row 12 (C12) which previously had no value
$ws.Range("C12").Value = 3

# Add a new row 20 with a note and reviewers
$ws.Range("B20").Value = "For Merging into Master Branch reviewers"
$ws.Range("E20").Value = "Mohamed Aboelazaiem"
$ws.Range("F20").Value = "Mahmoud Ahmed"

# Update the selection to match the final state
$ws.Range("F18").Select()
